# Adds columns I (I0) and J (IF) to the sheet, mirroring existing header /
# data-cell formatting conventions already used by columns B:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers "I0" and "IF" ---
# Copy the existing header formatting (bold font, border, centered) from H1
# so the new header cells reuse the workbook's existing style instead of
# creating new ones, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-70): numeric values for columns I and J ---
$data = @(
    @(9,9),
    @(6,6),
    @(8,8),
    @(8,8),
    @(6,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,8),
    @(8,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,10),
    @(7,7),
    @(7,8),
    @(7,7),
    @(9,9),
    @(8,8),
    @(4,5),
    @(6,7),
    @(8,8),
    @(6,6),
    @(8,8),
    @(9,9),
    @(7,7),
    @(6,6),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(8,8),
    @(10,11),
    @(9,9),
    @(9,9),
    @(4,4),
    @(10,11),
    @(6,6),
    @(7,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(8,8),
    @(11,11),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(9,9),
    @(4,5),
    @(6,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(6,6),
    @(5,5),
    @(6,6)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Output "Added columns I (I0) and J (IF) for rows 1-70"
